$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -3
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -10
$ws.Range("F11").Value = -1
$ws.Range("F13").Value = 0
$ws.Range("F16").Value = -9
